$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H), matching the formatting of the existing
# header cells (copy format from G1, which carries the bold/border/
# center-top style used by the other header cells).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for the new column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
